$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.785.15"
$ws.Range("E2").Value = "  -3.28%  "

$ws.Range("D3").Value = "'1.791.82"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.65%  "

$ws.Range("D5").Value = "'315.08"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").Value = "'0.5364"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").Value = "'0.3825"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").Value = "'0.07405"
$ws.Range("E9").Value = "  -1.87%  "

$ws.Range("D10").Value = "'41.66"
$ws.Range("E10").Value = "  -2.37%  "

$ws.Range("D11").Value = "'1.082"
$ws.Range("E11").Value = "  -3.42%  "

$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("D13").Value = "'6.171"
$ws.Range("E13").Value = "  -0.31%  "

$ws.Range("D14").Value = "'7.407"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "'20.22"
$ws.Range("E15").Value = "  -4.03%  "

$ws.Range("D16").Value = "'1.791.46"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "'88.23"
$ws.Range("E17").Value = "  -2.82%  "

$ws.Range("D18").Value = "'0.00001055"
$ws.Range("E18").Value = "  -1.17%  "

$ws.Range("D19").Value = "'0.06482"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "'17.24"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").Value = "'5.914"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("D23").Value = "'27.860.42"
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "  -1.33%  "

$ws.Range("D25").Value = "'2.090"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").Value = "'155.99"
$ws.Range("E26").Value = "  -3.03%  "

$ws.Range("D27").Value = "'20.22"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("D28").Value = "'2.002.90"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "'2.319"
$ws.Range("E29").Value = "  -2.77%  "

$ws.Range("D30").Value = "'121.51"
$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("D31").Value = "'0.1098"
$ws.Range("E31").Value = "  +5.44%  "

$ws.Range("D32").Value = "'1.102"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("D33").Value = "'3.652"
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("D34").Value = "'5.490"
$ws.Range("E34").Value = "  -3.58%  "

$ws.Range("D35").Value = "'0.06927"
$ws.Range("E35").Value = "  +6.41%  "

$ws.Range("D36").Value = "'0.2193"
$ws.Range("E36").Value = "  -3.21%  "

$ws.Range("D37").Value = "'0.02269"
$ws.Range("E37").Value = "  -2.28%  "

$ws.Range("D38").Value = "'5.031"
$ws.Range("E38").Value = "  -0.63%  "

$ws.Range("D39").Value = "'8.412"
$ws.Range("E39").Value = "  -6.00%  "

$ws.Range("D40").Value = "'11.29"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").Value = "'0.6076"
$ws.Range("E41").Value = "  -3.09%  "

$ws.Range("D42").Value = "'1.416"
$ws.Range("E42").Value = "  +1.70%  "

$ws.Range("D43").Value = "'1.159"
$ws.Range("E43").Value = "  -4.17%  "

$ws.Range("D44").Value = "'13.30"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").Value = "'3.674"
$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("D46").Value = "'0.5674"
$ws.Range("E46").Value = "  -3.79%  "

$ws.Range("D47").Value = "'123.93"
$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.900"
$ws.Range("E48").Value = "  -3.35%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.170"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("D50").Value = "'0.06784"
$ws.Range("E50").Value = "  -1.55%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000292"
$ws.Range("E51").Value = "  +35.23%  "
